$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels in row 2:
# "...מחוסנים עם דחף" (vaccinated with booster) -> "...מחוסנים" (vaccinated)
# "...מחוסנים ללא דחף" (vaccinated without booster) -> "...מחוסנים ללא תוקף" (vaccinated, no longer valid)
$ws.Range('B2').Value = 'חולים פעילים מחוסנים'
$ws.Range('C2').Value = 'חולים פעילים מחוסנים ללא תוקף'
$ws.Range('E2').Value = 'חולים פעילים מחוסנים ל-100 אלף תושבים'
$ws.Range('F2').Value = 'חולים פעילים מחוסנים ללא תוקף ל-100 אלף תושבים'
$ws.Range('H2').Value = 'חולים קשה מחוסנים'
$ws.Range('I2').Value = 'חולים קשה מחוסנים ללא תוקף'
$ws.Range('K2').Value = 'חולים קשה מחוסנים ל-100 אלף תושבים'
$ws.Range('L2').Value = 'חולים קשה מחוסנים ללא תוקף ל-100 אלף תושבים'

# Update data values for rows 3-12 (columns B:M) as of 2021-10-13
$data = New-Object 'object[,]' 10,12
$data[0,0] = 168
$data[0,1] = 0
$data[0,2] = 1411
$data[0,3] = 50.6
$data[0,4] = 0
$data[0,5] = 550.8
$data[0,6] = 0
$data[0,7] = 0
$data[0,8] = 1
$data[0,9] = 0
$data[0,10] = 0
$data[0,11] = 0.4
$data[1,0] = 192
$data[1,1] = 291
$data[1,2] = 684
$data[1,3] = 58.9
$data[1,4] = 289.4
$data[1,5] = 536.5
$data[1,6] = 0
$data[1,7] = 0
$data[1,8] = 1
$data[1,9] = 0
$data[1,10] = 0
$data[1,11] = 0.8
$data[2,0] = 481
$data[2,1] = 774
$data[2,2] = 1239
$data[2,3] = 63.5
$data[2,4] = 304.5
$data[2,5] = 544.5
$data[2,6] = 1
$data[2,7] = 3
$data[2,8] = 13
$data[2,9] = 0.1
$data[2,10] = 1.2
$data[2,11] = 5.7
$data[3,0] = 511
$data[3,1] = 846
$data[3,2] = 1210
$data[3,3] = 67.6
$data[3,4] = 374.5
$data[3,5] = 739.3
$data[3,6] = 0
$data[3,7] = 1
$data[3,8] = 22
$data[3,9] = 0
$data[3,10] = 0.4
$data[3,11] = 13.4
$data[4,0] = 384
$data[4,1] = 576
$data[4,2] = 843
$data[4,3] = 49.2
$data[4,4] = 356.6
$data[4,5] = 669.6
$data[4,6] = 2
$data[4,7] = 2
$data[4,8] = 47
$data[4,9] = 0.3
$data[4,10] = 1.2
$data[4,11] = 37.3
$data[5,0] = 231
$data[5,1] = 319
$data[5,2] = 518
$data[5,3] = 34.9
$data[5,4] = 318.3
$data[5,5] = 598.3
$data[5,6] = 8
$data[5,7] = 9
$data[5,8] = 70
$data[5,9] = 1.2
$data[5,10] = 9
$data[5,11] = 80.8
$data[6,0] = 209
$data[6,1] = 174
$data[6,2] = 337
$data[6,3] = 33.8
$data[6,4] = 320.1
$data[6,5] = 549.3
$data[6,6] = 18
$data[6,7] = 10
$data[6,8] = 62
$data[6,9] = 2.9
$data[6,10] = 18.4
$data[6,11] = 101.1
$data[7,0] = 157
$data[7,1] = 93
$data[7,2] = 230
$data[7,3] = 34.9
$data[7,4] = 318
$data[7,5] = 548.8
$data[7,6] = 16
$data[7,7] = 9
$data[7,8] = 47
$data[7,9] = 3.6
$data[7,10] = 30.8
$data[7,11] = 112.1
$data[8,0] = 109
$data[8,1] = 47
$data[8,2] = 128
$data[8,3] = 55.6
$data[8,4] = 260.9
$data[8,5] = 626.4
$data[8,6] = 10
$data[8,7] = 4
$data[8,8] = 23
$data[8,9] = 5.1
$data[8,10] = 22.2
$data[8,11] = 112.6
$data[9,0] = 35
$data[9,1] = 7
$data[9,2] = 32
$data[9,3] = 84.1
$data[9,4] = 94.6
$data[9,5] = 405
$data[9,6] = 2
$data[9,7] = 0
$data[9,8] = 10
$data[9,9] = 4.8
$data[9,10] = 0
$data[9,11] = 126.6

$ws.Range("B3:M12").Value = $data

